# Scheduled refresh: update FFXIV leve-crafting profit figures (currentAveragePrice*,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) per sheet/row, pulling the latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 3153.5344
$ws.Range("I137").Value = 932
$ws.Range("K137").Value = 2796
$ws.Range("M137").Value = -246
$ws = $wb.Worksheets.Item("ARM")
# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2756.875
$ws.Range("J63").Value = 2133.3333
$ws.Range("L63").Value = 2133.3333
$ws.Range("N63").Value = -3505.3333
# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2756.875
$ws.Range("J66").Value = 2133.3333
$ws.Range("L66").Value = 10666.6665
$ws.Range("N66").Value = -17530.6665
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 2839.7307
$ws.Range("I74").Value = 758.65
$ws.Range("J74").Value = 9776.666999999999
$ws.Range("K74").Value = 758.65
$ws.Range("L74").Value = 9776.666999999999
$ws.Range("M74").Value = 115.35
$ws.Range("N74").Value = -11524.667
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 2839.7307
$ws.Range("I77").Value = 758.65
$ws.Range("J77").Value = 9776.666999999999
$ws.Range("K77").Value = 3793.25
$ws.Range("L77").Value = 48883.335
$ws.Range("M77").Value = 574.75
$ws.Range("N77").Value = -57619.335
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 947.8095
$ws.Range("I97").Value = 693.9231
$ws.Range("J97").Value = 1360.375
$ws.Range("K97").Value = 693.9231
$ws.Range("L97").Value = 1360.375
$ws.Range("M97").Value = -197.9231
$ws.Range("N97").Value = -2352.375
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 6359.231
$ws.Range("I132").Value = 4367.448
$ws.Range("J132").Value = 12135.4
$ws.Range("K132").Value = 13102.344
$ws.Range("L132").Value = 36406.2
$ws.Range("M132").Value = -10572.344
$ws.Range("N132").Value = -41466.2
$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 42430.875
$ws.Range("I20").Value = 838.7692
$ws.Range("J20").Value = 91585.17999999999
$ws.Range("K20").Value = 838.7692
$ws.Range("L20").Value = 91585.17999999999
$ws.Range("M20").Value = -591.7692
$ws.Range("N20").Value = -92079.17999999999
# Row 80: Unbreaker / Titanium Ingot
$ws.Range("H80").Value = 471.7647
$ws.Range("I80").Value = 398.85715
$ws.Range("J80").Value = 522.8
$ws.Range("K80").Value = 398.85715
$ws.Range("L80").Value = 522.8
$ws.Range("M80").Value = 599.14285
$ws.Range("N80").Value = -2518.8
# Row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Range("H83").Value = 471.7647
$ws.Range("I83").Value = 398.85715
$ws.Range("J83").Value = 522.8
$ws.Range("K83").Value = 1994.28575
$ws.Range("L83").Value = 2614
$ws.Range("M83").Value = 2997.71425
$ws.Range("N83").Value = -12598
# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 2699.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2699.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2699.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4945.5
# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 2699.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2699.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 13497.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -24729.5
# Row 94: High Steal / High Steel Nugget
$ws.Range("H94").Value = 1730.5454
$ws.Range("I94").Value = 1612
$ws.Range("J94").Value = 2046.6666
$ws.Range("K94").Value = 1612
$ws.Range("L94").Value = 2046.6666
$ws.Range("M94").Value = -1161
$ws.Range("N94").Value = -2948.6666
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1346.4783
$ws.Range("I134").Value = 938.6875
$ws.Range("J134").Value = 2278.5715
$ws.Range("K134").Value = 2816.0625
$ws.Range("L134").Value = 6835.7145
$ws.Range("M134").Value = -281.0625
$ws.Range("N134").Value = -11905.7145
$ws = $wb.Worksheets.Item("CRP")
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 20004164
$ws.Range("I132").Value = 28576034
$ws.Range("J132").Value = 3136.2666
$ws.Range("K132").Value = 85728102
$ws.Range("L132").Value = 9408.799800000001
$ws.Range("M132").Value = -85725572
$ws.Range("N132").Value = -14468.7998
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 782
$ws.Range("I5").Value = 503.14285
$ws.Range("J5").Value = 998.8889
$ws.Range("K5").Value = 1509.42855
$ws.Range("L5").Value = 2996.6667
$ws.Range("M5").Value = -1397.42855
$ws.Range("N5").Value = -3220.6667
# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 871.65
$ws.Range("I122").Value = 551
$ws.Range("J122").Value = 1352.625
$ws.Range("K122").Value = 4959
$ws.Range("L122").Value = 12173.625
$ws.Range("M122").Value = -2509
$ws.Range("N122").Value = -17073.625
# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 782
$ws.Range("I135").Value = 503.14285
$ws.Range("J135").Value = 998.8889
$ws.Range("K135").Value = 4528.28565
$ws.Range("L135").Value = 8990.000100000001
$ws.Range("M135").Value = -1993.28565
$ws.Range("N135").Value = -14060.0001
$ws = $wb.Worksheets.Item("GSM")
# Row 123: Workplace Workout / Ametrine Ring of Fending
$ws.Range("H123").Value = 20244.375
$ws.Range("J123").Value = 20244.375
$ws.Range("L123").Value = 20244.375
$ws.Range("N123").Value = -25144.375
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 901.5
$ws.Range("I126").Value = 902.1111
$ws.Range("J126").Value = 899.6667
$ws.Range("K126").Value = 2706.3333
$ws.Range("L126").Value = 2699.0001
$ws.Range("M126").Value = -236.3332999999998
$ws.Range("N126").Value = -7639.0001
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 6795.5654
$ws.Range("I132").Value = 12873.6
$ws.Range("J132").Value = 2120.1538
$ws.Range("K132").Value = 38620.8
$ws.Range("L132").Value = 6360.4614
$ws.Range("M132").Value = -36090.8
$ws.Range("N132").Value = -11420.4614
$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 9600.823
$ws.Range("I132").Value = 15023.556
$ws.Range("J132").Value = 3500.25
$ws.Range("K132").Value = 45070.66800000001
$ws.Range("L132").Value = 10500.75
$ws.Range("M132").Value = -42540.66800000001
$ws.Range("N132").Value = -15560.75
$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 7408617.5
$ws.Range("I122").Value = 13334265
$ws.Range("J122").Value = 1559.0834
$ws.Range("K122").Value = 40002795
$ws.Range("L122").Value = 4677.2502
$ws.Range("M122").Value = -40000345
$ws.Range("N122").Value = -9577.2502
# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 679.7143
$ws.Range("I126").Value = 679.7143
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2039.1429
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = 430.8571000000002
$ws.Range("N126").ClearContents()
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 3610.5
$ws.Range("I136").Value = 3742.1667
$ws.Range("J136").Value = 3347.1667
$ws.Range("K136").Value = 11226.5001
$ws.Range("L136").Value = 10041.5001
$ws.Range("M136").Value = -8676.500100000001
$ws.Range("N136").Value = -15141.5001
